# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 23:05"

# --- Swap country labels: Chile/Pakistan (rows 22/23) ---
# Before: A22=Chile, A23=Pakistan -> After: A22=Pakistan, A23=Chile
$ws.Range("A22").Value = "Pakistan"
$ws.Range("A23").Value = "Chile"

# --- Swap country labels: Nueva Caledonia/Belice (rows 193/194) ---
# Before: A193=Nueva Caledonia, A194=Belice -> After: A193=Belice, A194=Nueva Caledonia
$ws.Range("A193").Value = "Belice"
$ws.Range("A194").Value = "Nueva Caledonia"

# --- Update numeric statistics cells ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1425890
$ws.Range("C4").Value = 17254
$ws.Range("E4").Value = 1034675
$ws.Range("G4").Value = 1492
$ws.Range("H4").Value = 84917

# Row 11 - Alemania
$ws.Range("B11").Value = 174098
$ws.Range("C11").Value = 927
$ws.Range("E11").Value = 17539
$ws.Range("G11").Value = 121
$ws.Range("H11").Value = 7859

# Row 17 - Canada
$ws.Range("B17").Value = 72200
$ws.Range("C17").Value = 1043
$ws.Range("D17").Value = 34927
$ws.Range("E17").Value = 31972

# Row 22 - Pakistan (after label swap)
$ws.Range("B22").Value = 35298
$ws.Range("C22").Value = 2624
$ws.Range("D22").Value = 8899
$ws.Range("E22").Value = 25638
$ws.Range("F22").Value = 111
$ws.Range("G22").Value = 37
$ws.Range("H22").Value = 761

# Row 23 - Chile (after label swap)
$ws.Range("B23").Value = 34381
$ws.Range("C23").Value = 2660
$ws.Range("D23").Value = 14865
$ws.Range("E23").Value = 19170
$ws.Range("F23").Value = 494
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 346

# Row 42 - Sudafrica
$ws.Range("B42").Value = 12074
$ws.Range("C42").Value = 724
$ws.Range("E42").Value = 7110

# Row 59 - Barein
$ws.Range("D59").Value = 2205
$ws.Range("E59").Value = 3601

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 2612
$ws.Range("C75").Value = 93
$ws.Range("D75").Value = 2076
$ws.Range("E75").Value = 525

# Row 82 - Costa de Marfil
$ws.Range("B82").Value = 1912
$ws.Range("C82").Value = 55
$ws.Range("D82").Value = 902
$ws.Range("E82").Value = 986
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 24

# Row 101 - Maldivas
$ws.Range("D101").Value = 40
$ws.Range("E101").Value = 911

# Row 108 - Niger
$ws.Range("B108").Value = 860
$ws.Range("C108").Value = 6
$ws.Range("D108").Value = 658
$ws.Range("E108").Value = 153
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 49

# Row 112 - Burkina Faso
$ws.Range("B112").Value = 773
$ws.Range("C112").Value = 7
$ws.Range("D112").Value = 592
$ws.Range("E112").Value = 130

# Row 167 - Islas Caimanes
$ws.Range("D167").Value = 54
$ws.Range("E167").Value = 31

# Row 193 - Belice (after label swap)
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Row 194 - Nueva Caledonia (after label swap)
$ws.Range("D194").Value = 18
$ws.Range("H194").Value = 0
